# Scheduled-runner market data refresh: updates cached market-board price /
# profit figures (columns H:N) on a handful of rows across several job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). A few rows also drop a
# trailing LeveProfitHQ (N) value entirely when the HQ price data is no
# longer available.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2764.2307
$ws.Range("I113").Value = 2101
$ws.Range("J113").Value = 3178.75
$ws.Range("K113").Value = 2101
$ws.Range("L113").Value = 3178.75
$ws.Range("M113").Value = 1153
$ws.Range("N113").Value = -9686.75

# Row 132
$ws.Range("H132").Value = 1595.0392
$ws.Range("I132").Value = 1672.8478
$ws.Range("J132").Value = 879.2
$ws.Range("K132").Value = 5018.5434
$ws.Range("L132").Value = 2637.6
$ws.Range("M132").Value = -2488.5434
$ws.Range("N132").Value = -7697.6

# Row 138
$ws.Range("H138").Value = 3404.4375
$ws.Range("I138").Value = 1614.8889
$ws.Range("J138").Value = 4710.324
$ws.Range("K138").Value = 4844.6667
$ws.Range("L138").Value = 14130.972
$ws.Range("M138").Value = 295.3333000000002
$ws.Range("N138").Value = -24410.972

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4986.2153
$ws.Range("I32").Value = 3415.0186
$ws.Range("J32").Value = 12699.363
$ws.Range("K32").Value = 3415.0186
$ws.Range("L32").Value = 12699.363
$ws.Range("M32").Value = -3128.0186
$ws.Range("N32").Value = -13273.363

# Row 45
$ws.Range("H45").Value = 1195.8846
$ws.Range("I45").Value = 1175.2667
$ws.Range("J45").Value = 1328.4286
$ws.Range("K45").Value = 1175.2667
$ws.Range("L45").Value = 1328.4286
$ws.Range("M45").Value = -798.2666999999999
$ws.Range("N45").Value = -2082.4286

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 110
$ws.Range("H110").Value = 970.96875
$ws.Range("I110").Value = 968.52
$ws.Range("J110").Value = 979.7143
$ws.Range("K110").Value = 968.52
$ws.Range("L110").Value = 979.7143
$ws.Range("M110").Value = 1076.48
$ws.Range("N110").Value = -5069.7143

$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 68395
$ws.Range("J13").Value = 68395
$ws.Range("L13").Value = 68395
$ws.Range("N13").Value = -68731

# Row 99
$ws.Range("H99").Value = 1299.8
$ws.Range("I99").Value = 1299.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1299.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 198.2
$ws.Range("N99").ClearContents()

# Row 105
$ws.Range("H105").Value = 1740769
$ws.Range("I105").Value = 3911286.2
$ws.Range("J105").Value = 4355.1
$ws.Range("K105").Value = 3911286.2
$ws.Range("L105").Value = 4355.1
$ws.Range("M105").Value = -3909539.2
$ws.Range("N105").Value = -7849.1

# Row 107
$ws.Range("H107").Value = 1286.5
$ws.Range("I107").Value = 1044.3334
$ws.Range("K107").Value = 1044.3334
$ws.Range("M107").Value = 875.6666

# Row 115
$ws.Range("H115").Value = 41484
$ws.Range("J115").Value = 41484
$ws.Range("L115").Value = 41484
$ws.Range("N115").Value = -44618

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1211.875
$ws.Range("I16").Value = 747.875
$ws.Range("J16").Value = 1443.875
$ws.Range("K16").Value = 747.875
$ws.Range("L16").Value = 1443.875
$ws.Range("M16").Value = -460.875
$ws.Range("N16").Value = -2017.875

# Row 31
$ws.Range("H31").Value = 5246.6665
$ws.Range("I31").Value = 5589.52
$ws.Range("J31").Value = 4467.4546
$ws.Range("K31").Value = 5589.52
$ws.Range("L31").Value = 4467.4546
$ws.Range("M31").Value = -5294.52
$ws.Range("N31").Value = -5057.4546

# Row 34
$ws.Range("H34").Value = 5246.6665
$ws.Range("I34").Value = 5589.52
$ws.Range("J34").Value = 4467.4546
$ws.Range("K34").Value = 5589.52
$ws.Range("L34").Value = 4467.4546
$ws.Range("M34").Value = -5387.52
$ws.Range("N34").Value = -4871.4546

# Row 58
$ws.Range("H58").Value = 2600410.2
$ws.Range("I58").Value = 5350141
$ws.Range("J58").Value = 3442
$ws.Range("K58").Value = 5350141
$ws.Range("L58").Value = 3442
$ws.Range("M58").Value = -5349938
$ws.Range("N58").Value = -3848

# Row 107
$ws.Range("H107").Value = 943.9032
$ws.Range("I107").Value = 957.75
$ws.Range("J107").Value = 918.7273
$ws.Range("K107").Value = 957.75
$ws.Range("L107").Value = 918.7273
$ws.Range("M107").Value = 962.25
$ws.Range("N107").Value = -4758.7273

# Row 113
$ws.Range("H113").Value = 1211.875
$ws.Range("I113").Value = 747.875
$ws.Range("J113").Value = 1443.875
$ws.Range("K113").Value = 747.875
$ws.Range("L113").Value = 1443.875
$ws.Range("M113").Value = 1422.125
$ws.Range("N113").Value = -5783.875

# Row 132
$ws.Range("H132").Value = 4144.6807
$ws.Range("I132").Value = 4120.05
$ws.Range("J132").Value = 4285.4287
$ws.Range("K132").Value = 12360.15
$ws.Range("L132").Value = 12856.2861
$ws.Range("M132").Value = -9830.150000000001
$ws.Range("N132").Value = -17916.2861

# Row 136
$ws.Range("H136").Value = 2600410.2
$ws.Range("I136").Value = 5350141
$ws.Range("J136").Value = 3442
$ws.Range("K136").Value = 16050423
$ws.Range("L136").Value = 10326
$ws.Range("M136").Value = -16047873
$ws.Range("N136").Value = -15426

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 3441.6155
$ws.Range("I75").Value = 475
$ws.Range("J75").Value = 3981
$ws.Range("K75").Value = 1425
$ws.Range("L75").Value = 11943
$ws.Range("M75").Value = -427
$ws.Range("N75").Value = -13939

# Row 78
$ws.Range("H78").Value = 3441.6155
$ws.Range("I78").Value = 475
$ws.Range("J78").Value = 3981
$ws.Range("K78").Value = 4275
$ws.Range("L78").Value = 35829
$ws.Range("M78").Value = 717
$ws.Range("N78").Value = -45813

# Row 113
$ws.Range("H113").Value = 711.9403
$ws.Range("I113").Value = 708.0645
$ws.Range("J113").Value = 760
$ws.Range("K113").Value = 2124.1935
$ws.Range("L113").Value = 2280
$ws.Range("M113").Value = 45.80650000000014
$ws.Range("N113").Value = -6620

# Row 123
$ws.Range("H123").Value = 2243.9473
$ws.Range("I123").Value = 1300
$ws.Range("J123").Value = 2420.9375
$ws.Range("K123").Value = 3900
$ws.Range("L123").Value = 7262.8125
$ws.Range("M123").Value = -1450
$ws.Range("N123").Value = -12162.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 440.35294
$ws.Range("J107").Value = 980.4
$ws.Range("L107").Value = 980.4
$ws.Range("N107").Value = -4820.4

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 640261.5600000001
$ws.Range("I61").Value = 21099.092
$ws.Range("J61").Value = 2002419
$ws.Range("K61").Value = 21099.092
$ws.Range("L61").Value = 2002419
$ws.Range("M61").Value = -20897.092
$ws.Range("N61").Value = -2002823

# Row 113
$ws.Range("H113").Value = 640261.5600000001
$ws.Range("I113").Value = 21099.092
$ws.Range("J113").Value = 2002419
$ws.Range("K113").Value = 21099.092
$ws.Range("L113").Value = 2002419
$ws.Range("M113").Value = -18929.092
$ws.Range("N113").Value = -2006759

# Row 136
$ws.Range("H136").Value = 4862.156
$ws.Range("I136").Value = 3051.682
$ws.Range("J136").Value = 6593.913
$ws.Range("K136").Value = 9155.045999999998
$ws.Range("L136").Value = 19781.739
$ws.Range("M136").Value = -6605.045999999998
$ws.Range("N136").Value = -24881.739

$ws = $wb.Worksheets.Item("WVR")
# Row 118
$ws.Range("H118").Value = 67600
$ws.Range("J118").Value = 67600
$ws.Range("L118").Value = 67600
$ws.Range("N118").Value = -70914

# Row 136
$ws.Range("H136").Value = 4280.08
$ws.Range("I136").Value = 4075.1353
$ws.Range("J136").Value = 4863.385
$ws.Range("K136").Value = 12225.4059
$ws.Range("L136").Value = 14590.155
$ws.Range("M136").Value = -9675.4059
$ws.Range("N136").Value = -19690.155
